$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 35, shifting existing rows 35:130 down to 36:131
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new week's data
$ws.Cells.Item(35, 1).Value = 4
$ws.Cells.Item(35, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(35, 3).Value = "Los Lagos"
$ws.Cells.Item(35, 4).Value = 44497
$ws.Cells.Item(35, 5).Value = 10
$ws.Cells.Item(35, 6).Value = 100112028
$ws.Cells.Item(35, 7).Value = "Sandia"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 300
$ws.Cells.Item(35, 11).Value = 900
$ws.Cells.Item(35, 12).Value = 900
$ws.Cells.Item(35, 13).Value = 900
$ws.Cells.Item(35, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(35, 15).Value = "Perú"
$ws.Cells.Item(35, 16).Value = 900
$ws.Cells.Item(35, 17).Value = 1
$ws.Cells.Item(35, 18).Value = "Hortaliza"
